$wb = $excel.ActiveWorkbook

# --- 1) Metadata sheet: update Title and Date values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B5").Value = 'Mapping Métier/CDA/FHIR : "Traitement déjà administré"'
$meta.Range("B8").Value = "2026-01-22T09:24:45+00:00"

# --- 2) "Mapping Table 1" sheet: insert a new root-equivalence row
#        (FRCDATraitement -> FRMedicationAdministrationDocument) right
#        after the URL row, pushing the existing field mappings down. ---
$map1 = $wb.Worksheets.Item("Mapping Table 1")

# Insert a new blank row at row 3 (existing rows 3..19 shift to 4..20)
$map1.Rows.Item(3).Insert()

# Copy formatting from the row below (now row 4) so the new row matches
# the rest of the table's style (borders / wrap / alignment).
$map1.Range("A4:E4").Copy()
$map1.Range("A3:E3").PasteSpecial(-4122)

# Fill in the new row's values
$map1.Range("A3").Value = "FRCDATraitement"
$map1.Range("C3").Value = "equivalent"
$map1.Range("D3").Value = "FRMedicationAdministrationDocument"
